$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted above the current row 238,
# pushing the existing rows 238-292 down to 239-293 (the data set grows
# by one record, as the commit message "Fruta / hortaliza, semanal" implies).
$ws.Rows(238).Insert()

# Populate the newly inserted row with the new observation. The
# constant descriptive columns (market, region, product taxonomy...)
# mirror the rest of the "Granada" block.
$ws.Range("A238").Value = 10
$ws.Range("B238").Value = "Vega Modelo de Temuco"
$ws.Range("C238").Value = "La Araucanía"
$ws.Range("D238").Value = 45173
$ws.Range("E238").Value = 9
$ws.Range("F238").Value = "Fruta"
$ws.Range("G238").Value = 100104
$ws.Range("H238").Value = "Frutos de pepita"
$ws.Range("I238").Value = 100104001
$ws.Range("J238").Value = "Granada"
$ws.Range("K238").Value = "Wonderfull"
$ws.Range("L238").Value = "Primera"
$ws.Range("M238").Value = 300
$ws.Range("N238").Value = 16000
$ws.Range("O238").Value = 16000
$ws.Range("P238").Value = 16000
$ws.Range("Q238").Value = "$/bandeja 10 kilos granel"
$ws.Range("R238").Value = "Provincia de Limarí"
$ws.Range("S238").Value = 1600
$ws.Range("T238").Value = 10
